# Append: 2025-09-30 12:48 JST
# Update the "取得日時" (retrieved-at) timestamp in column A (rows 2-9)
# of the "ランサーズ" sheet from 2025-09-30 12:38:27 to 2025-09-30 12:48:23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-09-30 12:38:27"
$newTimestamp = "2025-09-30 12:48:23"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 9 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value2 = $newTimestamp
    }
}
